$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 54, shifting existing rows 54-82 down to 55-83.
$ws.Rows.Item(54).Insert()

# Populate the newly inserted row 54 with the new record's data.
$ws.Range("A54").Value = 1
$ws.Range("B54").Value = "Agrícola del Norte S.A. de Arica"
$ws.Range("C54").Value = "Arica y Parinacota"
$ws.Range("D54").Value = 44917
$ws.Range("E54").Value = 15
$ws.Range("F54").Value = 100112009
$ws.Range("G54").Value = "Acelga"
$ws.Range("H54").Value = "Sin especificar"
$ws.Range("I54").Value = "Segunda"
$ws.Range("J54").Value = 350
$ws.Range("K54").Value = 1500
$ws.Range("L54").Value = 2000
$ws.Range("M54").Value = 1786
$ws.Range("N54").Value = "$/atado 2,5 a 3 kilos"
$ws.Range("O54").Value = "Región de Arica y Parinacota"
$ws.Range("P54").Value = 595
$ws.Range("Q54").Value = 3
$ws.Range("R54").Value = "Hortaliza"
